$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Context"
$ws.Range("B1").Value = "Text Response"

$ws.Range("A2").Value = "привет"
$ws.Range("B2").Value = "приветствую"

$ws.Range("A3").Value = "здравствуй"
$ws.Range("B3").Value = "здравствуйте"

$ws.Range("A4").Value = "добрый день!"
$ws.Range("B4").Value = "как дела?"

$ws.Range("A5").Value = "доброе утро!"

$ws.Range("A6").Value = "добрый вечер!"

$ws.Range("A7").Value = "доброго времени суток!"

$ws.Range("A9").Value = "У меня проблема"
$ws.Range("B9").Value = "приносим свои извинения, с вами в ближайшее время свяжется наш оператор"

$ws.Range("A10").Value = "Вознкла проблема"
$ws.Range("B10").Value = "Ожидайте, с вами скоро свяжется администратор"

$ws.Range("A11").Value = "Столкнулась с трудностями"

$ws.Range("A12").Value = "Это что вообще такое"

$ws.Range("A13").Value = "НЕ РАБТАЕТ"

$ws.Range("A14").Value = "Помогите"

$ws.Range("A15").Value = "Подскажите, пожалуйста"

$ws.Range("A16").Value = "Возник вопрос"

$ws.Range("A18").Value = "Спасибо"
$ws.Range("B18").Value = "Обращайтесь, мы во всем постараемся помочь"

$ws.Range("A19").Value = "Благодарю за помощь"
$ws.Range("B19").Value = "Не за что, если еще возникнут трудности, мы всегда поможем "

$ws.Range("A21").Value = "Пока"
$ws.Range("B21").Value = "Досвидания"

$ws.Range("A22").Value = "Досвидания"
$ws.Range("B22").Value = "До скорой встречи"

$ws.Range("A23").Value = "Удачи"

$ws.Range("A24").Value = "Услышимся"

[void]$ws.Range("E22").Select()
